$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header text (new crime data week) ---
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# --- Update weekly crime statistics table (rows 14-30) ---
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G14").Value = 1
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("H14").Value = -100
$ws.Range("H14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J14").Value = 1
$ws.Range("J14").NumberFormat = "#,##0"
$ws.Range("K14").Value = 300
$ws.Range("K14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("F15").Value = 1
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -37.5
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -62.5
$ws.Range("I16").Value = 73
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = 4.285714285714
$ws.Range("L16").Value = 25.862068965517
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 250
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 73.333333333333
$ws.Range("I17").Value = 222
$ws.Range("J17").Value = 194
$ws.Range("K17").Value = 14.432989690721
$ws.Range("L17").Value = 48
$ws.Range("D15").Copy($ws.Range("C18"))
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 75
$ws.Range("J18").Value = 77
$ws.Range("K18").Value = -18.181818181818
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -57.142857142857
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -25.581395348837
$ws.Range("I19").Value = 380
$ws.Range("J19").Value = 378
$ws.Range("K19").Value = 0.5291005291
$ws.Range("L19").Value = 32.867132867132
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -42.105263157894
$ws.Range("I20").Value = 99
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = -1
$ws.Range("L20").Value = 125
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -34.782608695652
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 90
$ws.Range("H21").Value = -11.111111111111
$ws.Range("I21").Value = 851
$ws.Range("J21").Value = 830
$ws.Range("K21").Value = 2.530120481927
$ws.Range("L21").Value = 35.079365079365
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -7.692307692307
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = 18.181818181818
$ws.Range("I24").Value = 1303
$ws.Range("J24").Value = 1240
$ws.Range("K24").Value = 5.08064516129
$ws.Range("L24").Value = 58.131067961165
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = -3.030303030303
$ws.Range("I25").Value = 454
$ws.Range("J25").Value = 441
$ws.Range("K25").Value = 2.947845804988
$ws.Range("L25").Value = 24.043715846994
$ws.Range("D15").Copy($ws.Range("C26"))
$ws.Range("D15").Copy($ws.Range("D26"))
$ws.Range("E15").Copy($ws.Range("E26"))
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 20
$ws.Range("K26").Value = 17.647058823529
$ws.Range("L26").Value = -28.571428571428
$ws.Range("F27").Value = 7
$ws.Range("I27").Value = 58
$ws.Range("K27").Value = 52.631578947368
$ws.Range("L27").Value = 38.095238095238
$ws.Range("D15").Copy($ws.Range("C30"))
$ws.Range("F30").Value = 3
$ws.Range("I30").Value = 8
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 166.666666666667
